$wb = $excel.ActiveWorkbook

# Rename "SB" -> "SIB"
$sib = $wb.Worksheets.Item("SB")
$sib.Name = "SIB"

# Add a new sheet "SIC" and move it to the end (after "Habitat Y")
$new = $wb.Worksheets.Add()
$new.Name = "SIC"
$new.Move($null, $wb.Worksheets.Item("Habitat Y"))

# Re-fetch by name (Move() can invalidate the old object handle) and make it the active sheet
$sic = $wb.Worksheets.Item("SIC")
$sic.Activate()
